# Automatic update of files.
# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from 2023-10-25 (45224) to 2023-11-03 (45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).Date.AddDays(45233)

$ws.Range("C2:C6").Value = $newDate
